$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (row 1 stays, rows 2-16 hold the new measurement series)
$bValues = @(
    0,
    42.72001872658765,
    45.27692569068709,
    35.35533905932738,
    40,
    35,
    46.09772228646444,
    47.43416490252569,
    36.40054944640259,
    31.6227766016838,
    44.7213595499958,
    43.01162633521314,
    47.16990566028302,
    30,
    41.23105625617661,
    32.01562118716424
)

# Column A (rows 2-16 hold the index 0..14)
$aValues = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
}

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Extend the bold/centered/bordered style already used by A2 down through A16
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A1").Select() | Out-Null
